$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update election results for row 2 (AVEIRO / OLIVEIRA DO BAIRRO)
$ws.Range("H2").Value = 224
$ws.Range("I2").Value = 651
$ws.Range("J2").Value = 2603
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 674
$ws.Range("N2").Value = 434
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 55
$ws.Range("S2").Value = 283
$ws.Range("T2").Value = 443
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 4072
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 3939
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 35
